$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 17 with data, mirroring the structure of previous rows
# Copy the style of the cell above (A16) so the date format (style 1) is reused
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(17, 1).Value = 42625.884247685186
$ws.Cells.Item(17, 2).Value = 33
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 13).Value = 0
$ws.Cells.Item(17, 14).Value = "Random"
